# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  - linked from the slide master (the deck's visible design,
#                            originally the "Integral" color scheme)
#   ppt/theme/theme2.xml  - linked from the notes master (originally the stock
#                            "Office Theme" color scheme)
#
# The recorded edit swaps the two themes' colour schemes: theme1.xml (the slide
# master's theme) ends up carrying the "Office" palette that used to live in
# theme2.xml. (Font scheme / format scheme are already byte-identical between
# the two theme parts, so only the 12 theme colours actually change.)
#
# Helper: build the OLE_COLOR (0x00BBGGRR) integer that PowerPoint's
# ThemeColorScheme.Colors(i).RGB setter expects from plain R,G,B byte values,
# so the persisted <a:srgbClr val="RRGGBB"/> comes out correct.
function ToOleColor($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# The slide master's theme (ppt/theme/theme1.xml) is the only theme part this
# object model exposes for editing; reach it via the Design collection.
$tcs = $p.Designs.Item(1).SlideMaster.Theme.ThemeColorScheme

# Target values == the palette that used to live in theme2.xml ("Office Theme").
# ThemeColorScheme.Colors index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$tcs.Colors(1).RGB  = ToOleColor 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = ToOleColor 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = ToOleColor 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = ToOleColor 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = ToOleColor 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = ToOleColor 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = ToOleColor 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = ToOleColor 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = ToOleColor 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = ToOleColor 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = ToOleColor 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = ToOleColor 0x95 0x4F 0x72   # folHlink
